# JDBCExecutorTemplate.xlsx edit:
# Rename the bean-style accessor placeholders to snake_case column names
# (firstName -> first_name, lastName -> last_name, catchPhrase -> catch_phrase,
#  aManager -> is_a_manager), matching the move to HSQLDB-backed JDBC tests.
# Also flips which sheet/tab is active: "Query" becomes the active tab
# instead of "Prepared".
#
# NOTE: values use single-quoted PowerShell strings throughout because the
# cell text itself contains literal ${...} sequences, which would otherwise
# be treated as PowerShell variable/subexpression syntax inside double
# quotes. Embedded literal single quotes are escaped by doubling ('').

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Query"
$ws2 = $wb.Worksheets.Item(2)   # "Prepared"

# --- Sheet "Query" (row 2 holds the sample/formula row) ---
$ws1.Range("A2").Value = '<jt:forEach items="${jdbc.execQuery(''SELECT * FROM employee'')}" var="employee" >${employee.first_name}'
$ws1.Range("B2").Value = '${employee.last_name}'
$ws1.Range("C2").Value = '${employee.salary}'
$ws1.Range("D2").Value = '${employee.title}'
$ws1.Range("E2").Value = '${employee.manager}'
$ws1.Range("F2").Value = '${employee.catch_phrase}'
$ws1.Range("G2").Value = '${employee.is_a_manager}</jt:forEach>'

# --- Sheet "Prepared" (row 3 holds the sample/formula row) ---
$ws2.Range("A3").Value = '<jt:forEach items="${jdbc.execQuery(''SELECT * FROM employee WHERE title = ?'', titleSearch)}" var="employee" >${employee.first_name}'
$ws2.Range("B3").Value = '${employee.last_name}'
$ws2.Range("C3").Value = '${employee.salary}'
$ws2.Range("D3").Value = '${employee.title}'
$ws2.Range("E3").Value = '${employee.manager}'
$ws2.Range("F3").Value = '${employee.catch_phrase}'
$ws2.Range("G3").Value = '${employee.is_a_manager}</jt:forEach></jt:forEach>'

# --- Make "Query" the active tab/sheet (was "Prepared") ---
$ws1.Activate()
